$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98, pushing the existing rows 98-113 down to 99-114.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly record.
$ws.Range("A98").Value = 5
$ws.Range("B98").Value = "Macroferia Regional de Talca"
$ws.Range("C98").Value = "Maule"
$ws.Range("D98").Value = 44889
$ws.Range("E98").Value = 7
$ws.Range("F98").Value = 100112026
$ws.Range("G98").Value = "Haba"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 250
$ws.Range("K98").Value = 10000
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = 10000
$ws.Range("N98").Value = "$/saco 25 kilos"
$ws.Range("O98").Value = "Provincia de Talca"
$ws.Range("P98").Value = 400
$ws.Range("Q98").Value = 25
$ws.Range("R98").Value = "Hortaliza"
